$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.58598166666667
$ws.Range("H2").Value = 82.75794500000001
$ws.Range("I2").Value = 0.2704460545904799
$ws.Range("J2").Value = 0.2704460545904799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 47.57896333333333
$ws.Range("N2").Value = 142.73689
$ws.Range("O2").Value = 0.450188452948237
$ws.Range("P2").Value = 0.4501884529482371
$ws.Range("Q2").Value = 1312.512410232339
$ws.Range("R2").Value = 11812.61169209105
$ws.Range("S2").Value = 0.1217516909220426
$ws.Range("T2").Value = 0.1217516909220426
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.58598166666667
$ws.Range("H3").Value = 82.75794500000001
$ws.Range("I3").Value = 0.2704460545904799
$ws.Range("J3").Value = 0.2704460545904799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.829723666666666
$ws.Range("N3").Value = 29.489171
$ws.Range("O3").Value = 0.09300808131111737
$ws.Range("P3").Value = 0.09300808131111739
$ws.Range("Q3").Value = 271.1625768570661
$ws.Range("R3").Value = 2440.463191713595
$ws.Range("S3").Value = 0.02515366863562225
$ws.Range("T3").Value = 0.02515366863562224
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.58598166666667
$ws.Range("H4").Value = 82.75794500000001
$ws.Range("I4").Value = 0.2704460545904799
$ws.Range("J4").Value = 0.2704460545904799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.933664
$ws.Range("N4").Value = 32.800992
$ws.Range("O4").Value = 0.1034534789405002
$ws.Range("P4").Value = 0.1034534789405003
$ws.Range("Q4").Value = 301.6158546534934
$ws.Range("R4").Value = 2714.54269188144
$ws.Range("S4").Value = 0.0279785852131176
$ws.Range("T4").Value = 0.02797858521311759
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 27.58598166666667
$ws.Range("H5").Value = 82.75794500000001
$ws.Range("I5").Value = 0.2704460545904799
$ws.Range("J5").Value = 0.2704460545904799
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 37.34441866666666
$ws.Range("N5").Value = 112.033256
$ws.Range("O5").Value = 0.3533499868001453
$ws.Range("P5").Value = 0.3533499868001453
$ws.Range("Q5").Value = 1030.182448690991
$ws.Range("R5").Value = 9271.642038218921
$ws.Range("S5").Value = 0.09556210981969745
$ws.Range("T5").Value = 0.09556210981969745
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 26.23504533333333
$ws.Range("H6").Value = 78.705136
$ws.Range("I6").Value = 0.2572018131577233
$ws.Range("J6").Value = 0.2572018131577233
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 47.57896333333333
$ws.Range("N6").Value = 142.73689
$ws.Range("O6").Value = 0.450188452948237
$ws.Range("P6").Value = 0.4501884529482371
$ws.Range("Q6").Value = 1248.236259963004
$ws.Range("R6").Value = 11234.12633966704
$ws.Range("S6").Value = 0.115789286360957
$ws.Range("T6").Value = 0.115789286360957
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 26.23504533333333
$ws.Range("H7").Value = 78.705136
$ws.Range("I7").Value = 0.2572018131577233
$ws.Range("J7").Value = 0.2572018131577233
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.829723666666666
$ws.Range("N7").Value = 29.489171
$ws.Range("O7").Value = 0.09300808131111737
$ws.Range("P7").Value = 0.09300808131111739
$ws.Range("Q7").Value = 257.8832460091396
$ws.Range("R7").Value = 2320.949214082256
$ws.Range("S7").Value = 0.02392184715154035
$ws.Range("T7").Value = 0.02392184715154035
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 26.23504533333333
$ws.Range("H8").Value = 78.705136
$ws.Range("I8").Value = 0.2572018131577233
$ws.Range("J8").Value = 0.2572018131577233
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.933664
$ws.Range("N8").Value = 32.800992
$ws.Range("O8").Value = 0.1034534789405002
$ws.Range("P8").Value = 0.1034534789405003
$ws.Range("Q8").Value = 286.8451706994347
$ws.Range("R8").Value = 2581.606536294912
$ws.Range("S8").Value = 0.02660842236097101
$ws.Range("T8").Value = 0.02660842236097101
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 26.23504533333333
$ws.Range("H9").Value = 78.705136
$ws.Range("I9").Value = 0.2572018131577233
$ws.Range("J9").Value = 0.2572018131577233
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 37.34441866666666
$ws.Range("N9").Value = 112.033256
$ws.Range("O9").Value = 0.3533499868001453
$ws.Range("P9").Value = 0.3533499868001453
$ws.Range("Q9").Value = 979.7325166669794
$ws.Range("R9").Value = 8817.592650002814
$ws.Range("S9").Value = 0.09088225728425497
$ws.Range("T9").Value = 0.09088225728425497
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.31506333333334
$ws.Range("H10").Value = 87.94519000000001
$ws.Range("I10").Value = 0.2873975381543141
$ws.Range("J10").Value = 0.2873975381543141
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 47.57896333333333
$ws.Range("N10").Value = 142.73689
$ws.Range("O10").Value = 0.450188452948237
$ws.Range("P10").Value = 0.4501884529482371
$ws.Range("Q10").Value = 1394.780323451011
$ws.Range("R10").Value = 12553.0229110591
$ws.Range("S10").Value = 0.1293830530828226
$ws.Range("T10").Value = 0.1293830530828226
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 29.31506333333334
$ws.Range("H11").Value = 87.94519000000001
$ws.Range("I11").Value = 0.2873975381543141
$ws.Range("J11").Value = 0.2873975381543141
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.829723666666666
$ws.Range("N11").Value = 29.489171
$ws.Range("O11").Value = 0.09300808131111737
$ws.Range("P11").Value = 0.09300808131111739
$ws.Range("Q11").Value = 288.1589718374989
$ws.Range("R11").Value = 2593.43074653749
$ws.Range("S11").Value = 0.0267302935972714
$ws.Range("T11").Value = 0.02673029359727141
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 29.31506333333334
$ws.Range("H12").Value = 87.94519000000001
$ws.Range("I12").Value = 0.2873975381543141
$ws.Range("J12").Value = 0.2873975381543141
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.933664
$ws.Range("N12").Value = 32.800992
$ws.Range("O12").Value = 0.1034534789405002
$ws.Range("P12").Value = 0.1034534789405003
$ws.Range("Q12").Value = 320.5210526253867
$ws.Range("R12").Value = 2884.68947362848
$ws.Range("S12").Value = 0.02973227516099895
$ws.Range("T12").Value = 0.02973227516099896
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 29.31506333333334
$ws.Range("H13").Value = 87.94519000000001
$ws.Range("I13").Value = 0.2873975381543141
$ws.Range("J13").Value = 0.2873975381543141
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 37.34441866666666
$ws.Range("N13").Value = 112.033256
$ws.Range("O13").Value = 0.3533499868001453
$ws.Range("P13").Value = 0.3533499868001453
$ws.Range("Q13").Value = 1094.753998359849
$ws.Range("R13").Value = 9852.785985238641
$ws.Range("S13").Value = 0.1015519163132211
$ws.Range("T13").Value = 0.1015519163132212
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 18.86569966666666
$ws.Range("H14").Value = 56.59709899999999
$ws.Range("I14").Value = 0.1849545940974826
$ws.Range("J14").Value = 0.1849545940974826
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 47.57896333333333
$ws.Range("N14").Value = 142.73689
$ws.Range("O14").Value = 0.450188452948237
$ws.Range("P14").Value = 0.4501884529482371
$ws.Range("Q14").Value = 897.610432698012
$ws.Range("R14").Value = 8078.493894282108
$ws.Range("S14").Value = 0.08326442258241483
$ws.Range("T14").Value = 0.08326442258241484
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 18.86569966666666
$ws.Range("H15").Value = 56.59709899999999
$ws.Range("I15").Value = 0.1849545940974826
$ws.Range("J15").Value = 0.1849545940974826
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 9.829723666666666
$ws.Range("N15").Value = 29.489171
$ws.Range("O15").Value = 0.09300808131111737
$ws.Range("P15").Value = 0.09300808131111739
$ws.Range("Q15").Value = 185.4446145016587
$ws.Range("R15").Value = 1669.001530514929
$ws.Range("S15").Value = 0.01720227192668337
$ws.Range("T15").Value = 0.01720227192668338
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 18.86569966666666
$ws.Range("H16").Value = 56.59709899999999
$ws.Range("I16").Value = 0.1849545940974826
$ws.Range("J16").Value = 0.1849545940974826
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.933664
$ws.Range("N16").Value = 32.800992
$ws.Range("O16").Value = 0.1034534789405002
$ws.Range("P16").Value = 0.1034534789405003
$ws.Range("Q16").Value = 206.2712212802453
$ws.Range("R16").Value = 1856.440991522208
$ws.Range("S16").Value = 0.01913419620541269
$ws.Range("T16").Value = 0.01913419620541269
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 18.86569966666666
$ws.Range("H17").Value = 56.59709899999999
$ws.Range("I17").Value = 0.1849545940974826
$ws.Range("J17").Value = 0.1849545940974826
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 37.34441866666666
$ws.Range("N17").Value = 112.033256
$ws.Range("O17").Value = 0.3533499868001453
$ws.Range("P17").Value = 0.3533499868001453
$ws.Range("Q17").Value = 704.5285867915936
$ws.Range("R17").Value = 6340.757281124343
$ws.Range("S17").Value = 0.06535370338297172
$ws.Range("T17").Value = 0.06535370338297172
